$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 579 ("「くつろぎなさい」استكن..."), shifting all rows
# below it up by one.
$ws.Rows.Item(579).Delete()
